# Auto-generated edit script: update crypto price/volume columns per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.891.46'
$ws.Range("E2").Value = '  -0.46%  '
$ws.Range("D3").Value = '2.300.03'
$ws.Range("E3").Value = '  -0.27%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.56'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.61%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.23'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.66%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.510'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.54%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.503'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.50%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.67'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.49%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0797'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.81%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '18.27'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.54%  '
$ws.Range("E13").Value = '  +0.99%  '
$ws.Range("E14").Value = '  -1.44%  '
$ws.Range("D15").Value = '2.659.84'
$ws.Range("E15").Value = '  -0.21%  '
$ws.Range("D16").Value = '2.302.09'
$ws.Range("E16").Value = '  +0.10%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.783'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.76%  '
$ws.Range("D18").Value = '42.836.56'
$ws.Range("E18").Value = '  -0.41%  '
$ws.Range("E19").Value = '  -1.99%  '
$ws.Range("E20").Value = '  -0.48%  '
$ws.Range("E21").Value = '  -1.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.53'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.47'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.69%  '
$ws.Range("E24").Value = '  -1.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.48'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.72%  '
$ws.Range("E26").Value = '  +0.98%  '
$ws.Range("E27").Value = '  +0.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '25.48'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.34%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '166.73'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("E30").Value = '  +1.20%  '
$ws.Range("E31").Value = '  -0.86%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '32.96'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.78'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.02'
$ws.Range("D35").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.37'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -4.83%  '
$ws.Range("E37").Value = '  -0.32%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0693'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.55%  '
$ws.Range("E39").Value = '  -0.68%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.75'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.87%  '
$ws.Range("E41").Value = '  -1.23%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.72'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.96%  '
$ws.Range("D43").Value = '2.012.55'
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("E44").Value = '  -2.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.16'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.02'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '17.96'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +3.03%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.78'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.89%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.92'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +6.17%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.98'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.89%  '
$ws.Range("D51").Value = '2.527.56'
$ws.Range("E51").Value = '  -0.42%  '
